# TC02_Canine_Filter_Breed-AmerStaffd.xlsx
# Add a new "TabName" column in front, labeling each row's source tab
# (CasesTab / SamplesTab / new FilesTab row), and add a brand new
# FilesTab row with its own Neo4j query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Insert a new column at the front for the "TabName" labels.
#    This shifts the existing A/B/C/D columns to B/C/D/E and keeps
#    their content + widths intact.
# ---------------------------------------------------------------
$ws.Columns("A").Insert()

# ---------------------------------------------------------------
# 2) Insert a brand new row 4 (FilesTab) below the existing SamplesTab
#    row, before filling in any new data.
# ---------------------------------------------------------------
$ws.Rows("4").Insert()

# ---------------------------------------------------------------
# 3) Fill in the new "TabName" column values.
# ---------------------------------------------------------------
$ws.Range("A1").Value2 = "TabName"
$ws.Range("A2").Value2 = "CasesTab"
$ws.Range("A3").Value2 = "SamplesTab"
$ws.Range("A4").Value2 = "FilesTab"

# ---------------------------------------------------------------
# 4) Fill in the new FilesTab row (row 4), reusing the same
#    StatQuery / file name text used by the other rows.
# ---------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value2 = $filesQuery
$ws.Range("C4").Value2 = $ws.Range("C2").Value2
$ws.Range("D4").Value2 = $ws.Range("D2").Value2
$ws.Range("E4").Value2 = $ws.Range("E2").Value2

# ---------------------------------------------------------------
# 5) Formatting: wrap text on the (now) query-text column C, and on
#    the new FilesTab query cell B4, matching the other query cells.
# ---------------------------------------------------------------
$ws.Range("C1:C4").WrapText = $true
$ws.Range("B4").WrapText = $true

# ---------------------------------------------------------------
# 6) Column widths (values chosen to land as close as this engine's
#    pixel-rounded column-width storage allows to the target widths
#    of 10.90625 and 123.36328125 respectively).
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10.0
$ws.Columns("C").ColumnWidth = 122.5

# ---------------------------------------------------------------
# 7) Row heights (row2/row3 keep their original heights automatically;
#    set the new row4 height explicitly).
# ---------------------------------------------------------------
$ws.Rows("4").RowHeight = 246.5

# ---------------------------------------------------------------
# 8) Sheet view: scroll down so row 3 is at top, zoom to 70%, and
#    move the active selection to A4.
# ---------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.Zoom = 70
$ws.Range("A4").Select()
